{"js": "const body = context.document.body;\nconst pairs = [\n  [\"412\u00f78=\", \"638\u00f78=\"],\n  [\"120\u00f76=\", \"369\u00f74=\"],\n  [\"145\u00f72=\", \"747\u00f72=\"],\n  [\"262\u00f72=\", \"582\u00f72=\"],\n  [\"646\u00f74=\", \"561\u00f75=\"],\n  [\"615\u00f74=\", \"505\u00f76=\"],\n  [\"342\u00f74=\", \"544\u00f75=\"],\n  [\"906\u00f76=\", \"331\u00f72=\"],\n  [\"791\u00f78=\", \"793\u00f72=\"],\n  [\"310\u00f77=\", \"604\u00f75=\"],\n  [\"883\u00f76=\", \"623\u00f73=\"],\n  [\"449\u00f78=\", \"604\u00f79=\"],\n  [\"366\u00f75=\", \"849\u00f79=\"],\n  [\"409\u00f79=\", \"631\u00f74=\"],\n  [\"273\u00f78=\", \"984\u00f76=\"],\n  [\"771\u00f72=\", \"425\u00f77=\"],\n  [\"242\u00f78=\", \"543\u00f75=\"],\n  [\"665\u00f76=\", \"385\u00f77=\"],\n  [\"139\u00f79=\", \"421\u00f72=\"],\n  [\"752\u00f78=\", \"754\u00f76=\"],\n  [\"407\u00f78=\", \"265\u00f79=\"],\n  [\"857\u00f74=\", \"509\u00f73=\"],\n  [\"661\u00f76=\", \"486\u00f78=\"],\n  [\"726\u00f72=\", \"301\u00f74=\"],\n  [\"260\u00f78=\", \"910\u00f72=\"]\n];\n\nconst searchResults = pairs.map(([oldText, newText]) => {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  return { results, oldText, newText };\n});\n\nawait context.sync();\n\nfor (const { results, oldText, newText } of searchResults) {\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${oldText}\", found ${results.items.length}`);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"412\u00f78=\", \"638\u00f78=\"),\n    @(\"120\u00f76=\", \"369\u00f74=\"),\n    @(\"145\u00f72=\", \"747\u00f72=\"),\n    @(\"262\u00f72=\", \"582\u00f72=\"),\n    @(\"646\u00f74=\", \"561\u00f75=\"),\n    @(\"615\u00f74=\", \"505\u00f76=\"),\n    @(\"342\u00f74=\", \"544\u00f75=\"),\n    @(\"906\u00f76=\", \"331\u00f72=\"),\n    @(\"791\u00f78=\", \"793\u00f72=\"),\n    @(\"310\u00f77=\", \"604\u00f75=\"),\n    @(\"883\u00f76=\", \"623\u00f73=\"),\n    @(\"449\u00f78=\", \"604\u00f79=\"),\n    @(\"366\u00f75=\", \"849\u00f79=\"),\n    @(\"409\u00f79=\", \"631\u00f74=\"),\n    @(\"273\u00f78=\", \"984\u00f76=\"),\n    @(\"771\u00f72=\", \"425\u00f77=\"),\n    @(\"242\u00f78=\", \"543\u00f75=\"),\n    @(\"665\u00f76=\", \"385\u00f77=\"),\n    @(\"139\u00f79=\", \"421\u00f72=\"),\n    @(\"752\u00f78=\", \"754\u00f76=\"),\n    @(\"407\u00f78=\", \"265\u00f79=\"),\n    @(\"857\u00f74=\", \"509\u00f73=\"),\n    @(\"661\u00f76=\", \"486\u00f78=\"),\n    @(\"726\u00f72=\", \"301\u00f74=\"),\n    @(\"260\u00f78=\", \"910\u00f72=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
